$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New temperature readings appended to row 14 (columns AA through AV)
$values = @(
    81.906999999999996,
    81.922615384615398,
    81.938230769230799,
    81.9538461538462,
    81.969461538461502,
    81.985076923076903,
    82.000692307692304,
    82.016307692307706,
    82.031923076923107,
    82.047538461538494,
    82.063153846153895,
    82.078769230769197,
    82.094384615384598,
    82.11,
    82.125615384615401,
    82.141230769230802,
    82.156846153846203,
    82.172461538461505,
    82.188076923076906,
    82.203692307692293,
    82.219307692307694,
    82.234923076923096
)

# AA is column 27
$startCol = 27
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(14, $startCol + $i).Value = $values[$i]
}

# New row 15 with bearing number 14 in column A
$ws.Cells.Item(15, 1).Value = 14

# Update the active selection to match the edited workbook
$ws.Range("B15").Select()
